$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 <- content from original row 7
$ws.Range("A6").Value = 130872697
$ws.Range("B6").Value = 8451
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 106545
$ws.Range("F6").Value = "Mindre märgborre"
$ws.Range("G6").Value = "Tomicus minor"
$ws.Range("H6").Value = "(Hartig, 1834)"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q6").Value = 570626
$ws.Range("R6").Value = 6736614
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = "Dalarna"
$ws.Range("U6").Value = "Falun"
$ws.Range("V6").Value = "Dalarna"
$ws.Range("W6").Value = "Svärdsjö"
$ws.Range("Y6").Value = "2026-01-25"
$ws.Range("Z6").Value = ""
$ws.Range("AA6").Value = "2026-01-25"
$ws.Range("AB6").Value = ""
$ws.Range("AC6").Value = ""
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Erik Danielsson"
$ws.Range("AX6").Value = "Erik Danielsson"
$ws.Range("AY6").Value = ""

# Row 7 <- content from original row 6
$ws.Range("A7").Value = 130872713
$ws.Range("B7").Value = 79244
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q7").Value = 570831
$ws.Range("R7").Value = 6736787
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = "Dalarna"
$ws.Range("U7").Value = "Falun"
$ws.Range("V7").Value = "Dalarna"
$ws.Range("W7").Value = "Svärdsjö"
$ws.Range("Y7").Value = "2026-01-25"
$ws.Range("Z7").Value = ""
$ws.Range("AA7").Value = "2026-01-25"
$ws.Range("AB7").Value = ""
$ws.Range("AC7").Value = ""
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = ""
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Erik Danielsson"
$ws.Range("AX7").Value = "Erik Danielsson"
$ws.Range("AY7").Value = ""

# Row 12 <- content from original row 13
$ws.Range("A12").Value = 130872717
$ws.Range("B12").Value = 79244
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("P12").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q12").Value = 571254
$ws.Range("R12").Value = 6736578
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = "Dalarna"
$ws.Range("U12").Value = "Falun"
$ws.Range("V12").Value = "Dalarna"
$ws.Range("W12").Value = "Svärdsjö"
$ws.Range("Y12").Value = "2026-01-25"
$ws.Range("Z12").Value = ""
$ws.Range("AA12").Value = "2026-01-25"
$ws.Range("AB12").Value = ""
$ws.Range("AC12").Value = ""
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AF12").Value = ""
$ws.Range("AG12").Value = $false
$ws.Range("AT12").Value = ""
$ws.Range("AW12").Value = "Erik Danielsson"
$ws.Range("AX12").Value = "Erik Danielsson"
$ws.Range("AY12").Value = ""

# Row 13 <- content from original row 12
$ws.Range("A13").Value = 130872695
$ws.Range("B13").Value = 79002
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 228912
$ws.Range("F13").Value = "Mörk kolflarnlav"
$ws.Range("G13").Value = "Carbonicola myrmecina"
$ws.Range("H13").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("P13").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q13").Value = 570816
$ws.Range("R13").Value = 6736802
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = "Dalarna"
$ws.Range("U13").Value = "Falun"
$ws.Range("V13").Value = "Dalarna"
$ws.Range("W13").Value = "Svärdsjö"
$ws.Range("Y13").Value = "2026-01-25"
$ws.Range("Z13").Value = ""
$ws.Range("AA13").Value = "2026-01-25"
$ws.Range("AB13").Value = ""
$ws.Range("AC13").Value = ""
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AF13").Value = ""
$ws.Range("AG13").Value = $false
$ws.Range("AT13").Value = ""
$ws.Range("AW13").Value = "Erik Danielsson"
$ws.Range("AX13").Value = "Erik Danielsson"
$ws.Range("AY13").Value = ""

# Row 15 <- content from original row 16
$ws.Range("A15").Value = 130872715
$ws.Range("B15").Value = 79244
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("P15").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q15").Value = 571193
$ws.Range("R15").Value = 6736684
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = "Dalarna"
$ws.Range("U15").Value = "Falun"
$ws.Range("V15").Value = "Dalarna"
$ws.Range("W15").Value = "Svärdsjö"
$ws.Range("Y15").Value = "2026-01-25"
$ws.Range("Z15").Value = ""
$ws.Range("AA15").Value = "2026-01-25"
$ws.Range("AB15").Value = ""
$ws.Range("AC15").Value = ""
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AF15").Value = ""
$ws.Range("AG15").Value = $false
$ws.Range("AT15").Value = ""
$ws.Range("AW15").Value = "Erik Danielsson"
$ws.Range("AX15").Value = "Erik Danielsson"
$ws.Range("AY15").Value = ""

# Row 16 <- content from original row 15
$ws.Range("A16").Value = 130872718
$ws.Range("B16").Value = 79244
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = "Garnlav"
$ws.Range("G16").Value = "Alectoria sarmentosa"
$ws.Range("H16").Value = "(Ach.) Ach."
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("P16").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q16").Value = 571142
$ws.Range("R16").Value = 6736599
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = "Dalarna"
$ws.Range("U16").Value = "Falun"
$ws.Range("V16").Value = "Dalarna"
$ws.Range("W16").Value = "Svärdsjö"
$ws.Range("Y16").Value = "2026-01-25"
$ws.Range("Z16").Value = ""
$ws.Range("AA16").Value = "2026-01-25"
$ws.Range("AB16").Value = ""
$ws.Range("AC16").Value = ""
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AF16").Value = ""
$ws.Range("AG16").Value = $false
$ws.Range("AT16").Value = ""
$ws.Range("AW16").Value = "Erik Danielsson"
$ws.Range("AX16").Value = "Erik Danielsson"
$ws.Range("AY16").Value = ""

# Row 21 <- content from original row 22
$ws.Range("A21").Value = 130872725
$ws.Range("B21").Value = 5177
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 100526
$ws.Range("F21").Value = "Bronshjon"
$ws.Range("G21").Value = "Callidium coriaceum"
$ws.Range("H21").Value = "Paykull, 1800"
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = ""
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("P21").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q21").Value = 570869
$ws.Range("R21").Value = 6736590
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = "Dalarna"
$ws.Range("U21").Value = "Falun"
$ws.Range("V21").Value = "Dalarna"
$ws.Range("W21").Value = "Svärdsjö"
$ws.Range("Y21").Value = "2026-01-25"
$ws.Range("Z21").Value = ""
$ws.Range("AA21").Value = "2026-01-25"
$ws.Range("AB21").Value = ""
$ws.Range("AC21").Value = ""
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AF21").Value = ""
$ws.Range("AG21").Value = $false
$ws.Range("AT21").Value = ""
$ws.Range("AW21").Value = "Erik Danielsson"
$ws.Range("AX21").Value = "Erik Danielsson"
$ws.Range("AY21").Value = ""

# Row 22 <- content from original row 21
$ws.Range("A22").Value = 130872716
$ws.Range("B22").Value = 79244
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("I22").Value = ""
$ws.Range("J22").Value = ""
$ws.Range("K22").Value = ""
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""
$ws.Range("P22").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q22").Value = 570988
$ws.Range("R22").Value = 6736647
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = "Dalarna"
$ws.Range("U22").Value = "Falun"
$ws.Range("V22").Value = "Dalarna"
$ws.Range("W22").Value = "Svärdsjö"
$ws.Range("Y22").Value = "2026-01-25"
$ws.Range("Z22").Value = ""
$ws.Range("AA22").Value = "2026-01-25"
$ws.Range("AB22").Value = ""
$ws.Range("AC22").Value = ""
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AF22").Value = ""
$ws.Range("AG22").Value = $false
$ws.Range("AT22").Value = ""
$ws.Range("AW22").Value = "Erik Danielsson"
$ws.Range("AX22").Value = "Erik Danielsson"
$ws.Range("AY22").Value = ""

# Row 25 <- content from original row 26
$ws.Range("A25").Value = 130983063
$ws.Range("B25").Value = 8451
$ws.Range("D25").Value = "LC"
$ws.Range("E25").Value = 106545
$ws.Range("F25").Value = "Mindre märgborre"
$ws.Range("G25").Value = "Tomicus minor"
$ws.Range("H25").Value = "(Hartig, 1834)"
$ws.Range("I25").Value = ""
$ws.Range("J25").Value = ""
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = "äldre gnagspår"
$ws.Range("N25").Value = ""
$ws.Range("P25").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q25").Value = 570956
$ws.Range("R25").Value = 6736657
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = "Dalarna"
$ws.Range("U25").Value = "Falun"
$ws.Range("V25").Value = "Dalarna"
$ws.Range("W25").Value = "Svärdsjö"
$ws.Range("Y25").Value = "2026-01-31"
$ws.Range("Z25").Value = "09:32"
$ws.Range("AA25").Value = "2026-01-31"
$ws.Range("AB25").Value = "09:32"
$ws.Range("AC25").Value = ""
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AF25").Value = ""
$ws.Range("AG25").Value = $false
$ws.Range("AT25").Value = ""
$ws.Range("AW25").Value = "Bo karlstens"
$ws.Range("AX25").Value = "Bo karlstens"
$ws.Range("AY25").Value = ""

# Row 26 <- content from original row 25
$ws.Range("A26").Value = 130979083
$ws.Range("B26").Value = 57073
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 100138
$ws.Range("F26").Value = "Tjäder"
$ws.Range("G26").Value = "Tetrao urogallus"
$ws.Range("H26").Value = "Linnaeus, 1758"
$ws.Range("I26").Value = ""
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = ""
$ws.Range("L26").Value = ""
$ws.Range("M26").Value = ""
$ws.Range("N26").Value = ""
$ws.Range("P26").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q26").Value = 570745
$ws.Range("R26").Value = 6736794
$ws.Range("S26").Value = 1
$ws.Range("T26").Value = "Dalarna"
$ws.Range("U26").Value = "Falun"
$ws.Range("V26").Value = "Dalarna"
$ws.Range("W26").Value = "Svärdsjö"
$ws.Range("Y26").Value = "2026-01-31"
$ws.Range("Z26").Value = ""
$ws.Range("AA26").Value = "2026-01-31"
$ws.Range("AB26").Value = ""
$ws.Range("AC26").Value = "Färsk spillning"
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AF26").Value = ""
$ws.Range("AG26").Value = $false
$ws.Range("AT26").Value = ""
$ws.Range("AW26").Value = "Erik Danielsson"
$ws.Range("AX26").Value = "Erik Danielsson"
$ws.Range("AY26").Value = ""

# Row 28 <- content from original row 29
$ws.Range("A28").Value = 130979095
$ws.Range("B28").Value = 79244
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("I28").Value = ""
$ws.Range("J28").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("L28").Value = ""
$ws.Range("M28").Value = ""
$ws.Range("N28").Value = ""
$ws.Range("P28").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q28").Value = 571648
$ws.Range("R28").Value = 6736781
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = "Dalarna"
$ws.Range("U28").Value = "Falun"
$ws.Range("V28").Value = "Dalarna"
$ws.Range("W28").Value = "Svärdsjö"
$ws.Range("Y28").Value = "2026-01-31"
$ws.Range("Z28").Value = ""
$ws.Range("AA28").Value = "2026-01-31"
$ws.Range("AB28").Value = ""
$ws.Range("AC28").Value = ""
$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AF28").Value = ""
$ws.Range("AG28").Value = $false
$ws.Range("AT28").Value = ""
$ws.Range("AW28").Value = "Erik Danielsson"
$ws.Range("AX28").Value = "Erik Danielsson"
$ws.Range("AY28").Value = ""

# Row 29 <- content from original row 28
$ws.Range("A29").Value = 130979098
$ws.Range("B29").Value = 79244
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("I29").Value = ""
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = ""
$ws.Range("P29").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q29").Value = 571475
$ws.Range("R29").Value = 6736616
$ws.Range("S29").Value = 1
$ws.Range("T29").Value = "Dalarna"
$ws.Range("U29").Value = "Falun"
$ws.Range("V29").Value = "Dalarna"
$ws.Range("W29").Value = "Svärdsjö"
$ws.Range("Y29").Value = "2026-01-31"
$ws.Range("Z29").Value = ""
$ws.Range("AA29").Value = "2026-01-31"
$ws.Range("AB29").Value = ""
$ws.Range("AC29").Value = ""
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AF29").Value = ""
$ws.Range("AG29").Value = $false
$ws.Range("AT29").Value = ""
$ws.Range("AW29").Value = "Erik Danielsson"
$ws.Range("AX29").Value = "Erik Danielsson"
$ws.Range("AY29").Value = ""

# Row 30 <- content from original row 31
$ws.Range("A30").Value = 130983077
$ws.Range("B30").Value = 5177
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 100526
$ws.Range("F30").Value = "Bronshjon"
$ws.Range("G30").Value = "Callidium coriaceum"
$ws.Range("H30").Value = "Paykull, 1800"
$ws.Range("I30").Value = ""
$ws.Range("J30").Value = ""
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = "färska gnagspår"
$ws.Range("N30").Value = ""
$ws.Range("P30").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q30").Value = 571069
$ws.Range("R30").Value = 6736680
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = "Dalarna"
$ws.Range("U30").Value = "Falun"
$ws.Range("V30").Value = "Dalarna"
$ws.Range("W30").Value = "Svärdsjö"
$ws.Range("Y30").Value = "2026-01-31"
$ws.Range("Z30").Value = "11:22"
$ws.Range("AA30").Value = "2026-01-31"
$ws.Range("AB30").Value = "11:22"
$ws.Range("AC30").Value = ""
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AF30").Value = ""
$ws.Range("AG30").Value = $false
$ws.Range("AT30").Value = ""
$ws.Range("AW30").Value = "Bo karlstens"
$ws.Range("AX30").Value = "Bo karlstens"
$ws.Range("AY30").Value = ""

# Row 31 <- content from original row 30
$ws.Range("A31").Value = 130979089
$ws.Range("B31").Value = 57076
$ws.Range("D31").Value = "LC"
$ws.Range("E31").Value = 102613
$ws.Range("F31").Value = "Orre"
$ws.Range("G31").Value = "Lyrurus tetrix"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("I31").Value = ""
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = ""
$ws.Range("L31").Value = ""
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = ""
$ws.Range("P31").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q31").Value = 570598
$ws.Range("R31").Value = 6736697
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = "Dalarna"
$ws.Range("U31").Value = "Falun"
$ws.Range("V31").Value = "Dalarna"
$ws.Range("W31").Value = "Svärdsjö"
$ws.Range("Y31").Value = "2026-01-31"
$ws.Range("Z31").Value = ""
$ws.Range("AA31").Value = "2026-01-31"
$ws.Range("AB31").Value = ""
$ws.Range("AC31").Value = "Överflygande"
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AF31").Value = ""
$ws.Range("AG31").Value = $false
$ws.Range("AT31").Value = ""
$ws.Range("AW31").Value = "Erik Danielsson"
$ws.Range("AX31").Value = "Erik Danielsson"
$ws.Range("AY31").Value = ""

# Row 32 <- content from original row 33
$ws.Range("A32").Value = 130983618
$ws.Range("B32").Value = 79244
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = "Garnlav"
$ws.Range("G32").Value = "Alectoria sarmentosa"
$ws.Range("H32").Value = "(Ach.) Ach."
$ws.Range("I32").Value = ""
$ws.Range("J32").Value = ""
$ws.Range("K32").Value = ""
$ws.Range("L32").Value = ""
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = ""
$ws.Range("P32").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q32").Value = 570808
$ws.Range("R32").Value = 6736568
$ws.Range("S32").Value = 10
$ws.Range("T32").Value = "Dalarna"
$ws.Range("U32").Value = "Falun"
$ws.Range("V32").Value = "Dalarna"
$ws.Range("W32").Value = "Svärdsjö"
$ws.Range("Y32").Value = "2026-01-31"
$ws.Range("Z32").Value = "09:07"
$ws.Range("AA32").Value = "2026-01-31"
$ws.Range("AB32").Value = "09:07"
$ws.Range("AC32").Value = ""
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AF32").Value = ""
$ws.Range("AG32").Value = $false
$ws.Range("AT32").Value = ""
$ws.Range("AW32").Value = "Göran Ehn"
$ws.Range("AX32").Value = "Göran Ehn"
$ws.Range("AY32").Value = ""

# Row 33 <- content from original row 32
$ws.Range("A33").Value = 130983060
$ws.Range("B33").Value = 8451
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 106545
$ws.Range("F33").Value = "Mindre märgborre"
$ws.Range("G33").Value = "Tomicus minor"
$ws.Range("H33").Value = "(Hartig, 1834)"
$ws.Range("I33").Value = ""
$ws.Range("J33").Value = ""
$ws.Range("K33").Value = ""
$ws.Range("L33").Value = ""
$ws.Range("M33").Value = "äldre gnagspår"
$ws.Range("N33").Value = ""
$ws.Range("P33").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q33").Value = 570988
$ws.Range("R33").Value = 6736721
$ws.Range("S33").Value = 10
$ws.Range("T33").Value = "Dalarna"
$ws.Range("U33").Value = "Falun"
$ws.Range("V33").Value = "Dalarna"
$ws.Range("W33").Value = "Svärdsjö"
$ws.Range("Y33").Value = "2026-01-31"
$ws.Range("Z33").Value = "11:29"
$ws.Range("AA33").Value = "2026-01-31"
$ws.Range("AB33").Value = "11:29"
$ws.Range("AC33").Value = ""
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AF33").Value = ""
$ws.Range("AG33").Value = $false
$ws.Range("AT33").Value = ""
$ws.Range("AW33").Value = "Bo karlstens"
$ws.Range("AX33").Value = "Bo karlstens"
$ws.Range("AY33").Value = ""

# Row 48 <- content from original row 49
$ws.Range("A48").Value = 130983056
$ws.Range("B48").Value = 57884
$ws.Range("D48").Value = "NT"
$ws.Range("E48").Value = 100109
$ws.Range("F48").Value = "Tretåig hackspett"
$ws.Range("G48").Value = "Picoides tridactylus"
$ws.Range("H48").Value = "(Linnaeus, 1758)"
$ws.Range("I48").Value = ""
$ws.Range("J48").Value = ""
$ws.Range("K48").Value = ""
$ws.Range("L48").Value = ""
$ws.Range("M48").Value = "färska spår"
$ws.Range("N48").Value = ""
$ws.Range("P48").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q48").Value = 570823
$ws.Range("R48").Value = 6736624
$ws.Range("S48").Value = 10
$ws.Range("T48").Value = "Dalarna"
$ws.Range("U48").Value = "Falun"
$ws.Range("V48").Value = "Dalarna"
$ws.Range("W48").Value = "Svärdsjö"
$ws.Range("Y48").Value = "2026-01-31"
$ws.Range("Z48").Value = "09:12"
$ws.Range("AA48").Value = "2026-01-31"
$ws.Range("AB48").Value = "09:12"
$ws.Range("AC48").Value = ""
$ws.Range("AD48").Value = $false
$ws.Range("AE48").Value = $true
$ws.Range("AF48").Value = ""
$ws.Range("AG48").Value = $false
$ws.Range("AT48").Value = ""
$ws.Range("AW48").Value = "Bo karlstens"
$ws.Range("AX48").Value = "Bo karlstens"
$ws.Range("AY48").Value = ""

# Row 49 <- content from original row 50
$ws.Range("A49").Value = 130983061
$ws.Range("B49").Value = 8451
$ws.Range("D49").Value = "LC"
$ws.Range("E49").Value = 106545
$ws.Range("F49").Value = "Mindre märgborre"
$ws.Range("G49").Value = "Tomicus minor"
$ws.Range("H49").Value = "(Hartig, 1834)"
$ws.Range("I49").Value = ""
$ws.Range("J49").Value = ""
$ws.Range("K49").Value = ""
$ws.Range("L49").Value = ""
$ws.Range("M49").Value = "äldre gnagspår"
$ws.Range("N49").Value = ""
$ws.Range("P49").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q49").Value = 571287
$ws.Range("R49").Value = 6736565
$ws.Range("S49").Value = 10
$ws.Range("T49").Value = "Dalarna"
$ws.Range("U49").Value = "Falun"
$ws.Range("V49").Value = "Dalarna"
$ws.Range("W49").Value = "Svärdsjö"
$ws.Range("Y49").Value = "2026-01-31"
$ws.Range("Z49").Value = "10:34"
$ws.Range("AA49").Value = "2026-01-31"
$ws.Range("AB49").Value = "10:34"
$ws.Range("AC49").Value = ""
$ws.Range("AD49").Value = $false
$ws.Range("AE49").Value = $false
$ws.Range("AF49").Value = ""
$ws.Range("AG49").Value = $false
$ws.Range("AT49").Value = ""
$ws.Range("AW49").Value = "Bo karlstens"
$ws.Range("AX49").Value = "Bo karlstens"
$ws.Range("AY49").Value = ""

# Row 50 <- content from original row 48
$ws.Range("A50").Value = 130983067
$ws.Range("B50").Value = 79244
$ws.Range("D50").Value = "NT"
$ws.Range("E50").Value = 6425
$ws.Range("F50").Value = "Garnlav"
$ws.Range("G50").Value = "Alectoria sarmentosa"
$ws.Range("H50").Value = "(Ach.) Ach."
$ws.Range("I50").Value = ""
$ws.Range("J50").Value = ""
$ws.Range("K50").Value = ""
$ws.Range("L50").Value = ""
$ws.Range("M50").Value = ""
$ws.Range("N50").Value = ""
$ws.Range("P50").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q50").Value = 571052
$ws.Range("R50").Value = 6736642
$ws.Range("S50").Value = 10
$ws.Range("T50").Value = "Dalarna"
$ws.Range("U50").Value = "Falun"
$ws.Range("V50").Value = "Dalarna"
$ws.Range("W50").Value = "Svärdsjö"
$ws.Range("Y50").Value = "2026-01-31"
$ws.Range("Z50").Value = "09:38"
$ws.Range("AA50").Value = "2026-01-31"
$ws.Range("AB50").Value = "09:38"
$ws.Range("AC50").Value = ""
$ws.Range("AD50").Value = $false
$ws.Range("AE50").Value = $false
$ws.Range("AF50").Value = ""
$ws.Range("AG50").Value = $false
$ws.Range("AT50").Value = ""
$ws.Range("AW50").Value = "Bo karlstens"
$ws.Range("AX50").Value = "Bo karlstens"
$ws.Range("AY50").Value = ""

# Row 51 <- content from original row 53
$ws.Range("A51").Value = 130983612
$ws.Range("B51").Value = 79244
$ws.Range("D51").Value = "NT"
$ws.Range("E51").Value = 6425
$ws.Range("F51").Value = "Garnlav"
$ws.Range("G51").Value = "Alectoria sarmentosa"
$ws.Range("H51").Value = "(Ach.) Ach."
$ws.Range("I51").Value = ""
$ws.Range("J51").Value = ""
$ws.Range("K51").Value = ""
$ws.Range("L51").Value = ""
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = ""
$ws.Range("P51").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q51").Value = 571254
$ws.Range("R51").Value = 6736555
$ws.Range("S51").Value = 10
$ws.Range("T51").Value = "Dalarna"
$ws.Range("U51").Value = "Falun"
$ws.Range("V51").Value = "Dalarna"
$ws.Range("W51").Value = "Svärdsjö"
$ws.Range("Y51").Value = "2026-01-31"
$ws.Range("Z51").Value = "10:29"
$ws.Range("AA51").Value = "2026-01-31"
$ws.Range("AB51").Value = "10:29"
$ws.Range("AC51").Value = ""
$ws.Range("AD51").Value = $false
$ws.Range("AE51").Value = $false
$ws.Range("AF51").Value = ""
$ws.Range("AG51").Value = $false
$ws.Range("AT51").Value = ""
$ws.Range("AW51").Value = "Göran Ehn"
$ws.Range("AX51").Value = "Göran Ehn"
$ws.Range("AY51").Value = ""

# Row 52 <- content from original row 51
$ws.Range("A52").Value = 130983616
$ws.Range("B52").Value = 79244
$ws.Range("D52").Value = "NT"
$ws.Range("E52").Value = 6425
$ws.Range("F52").Value = "Garnlav"
$ws.Range("G52").Value = "Alectoria sarmentosa"
$ws.Range("H52").Value = "(Ach.) Ach."
$ws.Range("I52").Value = ""
$ws.Range("J52").Value = ""
$ws.Range("K52").Value = ""
$ws.Range("L52").Value = ""
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = ""
$ws.Range("P52").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q52").Value = 570925
$ws.Range("R52").Value = 6736674
$ws.Range("S52").Value = 10
$ws.Range("T52").Value = "Dalarna"
$ws.Range("U52").Value = "Falun"
$ws.Range("V52").Value = "Dalarna"
$ws.Range("W52").Value = "Svärdsjö"
$ws.Range("Y52").Value = "2026-01-31"
$ws.Range("Z52").Value = "09:16"
$ws.Range("AA52").Value = "2026-01-31"
$ws.Range("AB52").Value = "09:16"
$ws.Range("AC52").Value = ""
$ws.Range("AD52").Value = $false
$ws.Range("AE52").Value = $false
$ws.Range("AF52").Value = ""
$ws.Range("AG52").Value = $false
$ws.Range("AT52").Value = ""
$ws.Range("AW52").Value = "Göran Ehn"
$ws.Range("AX52").Value = "Göran Ehn"
$ws.Range("AY52").Value = ""

# Row 53 <- content from original row 52
$ws.Range("A53").Value = 130983607
$ws.Range("B53").Value = 79244
$ws.Range("D53").Value = "NT"
$ws.Range("E53").Value = 6425
$ws.Range("F53").Value = "Garnlav"
$ws.Range("G53").Value = "Alectoria sarmentosa"
$ws.Range("H53").Value = "(Ach.) Ach."
$ws.Range("I53").Value = ""
$ws.Range("J53").Value = ""
$ws.Range("K53").Value = ""
$ws.Range("L53").Value = ""
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = ""
$ws.Range("P53").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q53").Value = 571351
$ws.Range("R53").Value = 6736798
$ws.Range("S53").Value = 10
$ws.Range("T53").Value = "Dalarna"
$ws.Range("U53").Value = "Falun"
$ws.Range("V53").Value = "Dalarna"
$ws.Range("W53").Value = "Svärdsjö"
$ws.Range("Y53").Value = "2026-01-31"
$ws.Range("Z53").Value = "10:59"
$ws.Range("AA53").Value = "2026-01-31"
$ws.Range("AB53").Value = "10:59"
$ws.Range("AC53").Value = ""
$ws.Range("AD53").Value = $false
$ws.Range("AE53").Value = $false
$ws.Range("AF53").Value = ""
$ws.Range("AG53").Value = $false
$ws.Range("AT53").Value = ""
$ws.Range("AW53").Value = "Göran Ehn"
$ws.Range("AX53").Value = "Göran Ehn"
$ws.Range("AY53").Value = ""

# Row 54 <- content from original row 55
$ws.Range("A54").Value = 130983611
$ws.Range("B54").Value = 79244
$ws.Range("D54").Value = "NT"
$ws.Range("E54").Value = 6425
$ws.Range("F54").Value = "Garnlav"
$ws.Range("G54").Value = "Alectoria sarmentosa"
$ws.Range("H54").Value = "(Ach.) Ach."
$ws.Range("I54").Value = ""
$ws.Range("J54").Value = ""
$ws.Range("K54").Value = ""
$ws.Range("L54").Value = ""
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = ""
$ws.Range("P54").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q54").Value = 571283
$ws.Range("R54").Value = 6736557
$ws.Range("S54").Value = 10
$ws.Range("T54").Value = "Dalarna"
$ws.Range("U54").Value = "Falun"
$ws.Range("V54").Value = "Dalarna"
$ws.Range("W54").Value = "Svärdsjö"
$ws.Range("Y54").Value = "2026-01-31"
$ws.Range("Z54").Value = "10:33"
$ws.Range("AA54").Value = "2026-01-31"
$ws.Range("AB54").Value = "10:33"
$ws.Range("AC54").Value = ""
$ws.Range("AD54").Value = $false
$ws.Range("AE54").Value = $false
$ws.Range("AF54").Value = ""
$ws.Range("AG54").Value = $false
$ws.Range("AT54").Value = ""
$ws.Range("AW54").Value = "Göran Ehn"
$ws.Range("AX54").Value = "Göran Ehn"
$ws.Range("AY54").Value = ""

# Row 55 <- content from original row 54
$ws.Range("A55").Value = 130983062
$ws.Range("B55").Value = 8451
$ws.Range("D55").Value = "LC"
$ws.Range("E55").Value = 106545
$ws.Range("F55").Value = "Mindre märgborre"
$ws.Range("G55").Value = "Tomicus minor"
$ws.Range("H55").Value = "(Hartig, 1834)"
$ws.Range("I55").Value = ""
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""
$ws.Range("L55").Value = ""
$ws.Range("M55").Value = "äldre gnagspår"
$ws.Range("N55").Value = ""
$ws.Range("P55").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q55").Value = 571286
$ws.Range("R55").Value = 6736563
$ws.Range("S55").Value = 10
$ws.Range("T55").Value = "Dalarna"
$ws.Range("U55").Value = "Falun"
$ws.Range("V55").Value = "Dalarna"
$ws.Range("W55").Value = "Svärdsjö"
$ws.Range("Y55").Value = "2026-01-31"
$ws.Range("Z55").Value = "10:34"
$ws.Range("AA55").Value = "2026-01-31"
$ws.Range("AB55").Value = "10:34"
$ws.Range("AC55").Value = ""
$ws.Range("AD55").Value = $false
$ws.Range("AE55").Value = $false
$ws.Range("AF55").Value = ""
$ws.Range("AG55").Value = $false
$ws.Range("AT55").Value = ""
$ws.Range("AW55").Value = "Bo karlstens"
$ws.Range("AX55").Value = "Bo karlstens"
$ws.Range("AY55").Value = ""

# Row 56 <- content from original row 59
$ws.Range("A56").Value = 130979086
$ws.Range("B56").Value = 91830
$ws.Range("D56").Value = "NT"
$ws.Range("E56").Value = 5442
$ws.Range("F56").Value = "Tallticka"
$ws.Range("G56").Value = "Porodaedalea pini"
$ws.Range("H56").Value = "(Brot.) Murrill"
$ws.Range("I56").Value = ""
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""
$ws.Range("L56").Value = ""
$ws.Range("M56").Value = ""
$ws.Range("N56").Value = ""
$ws.Range("P56").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q56").Value = 571361
$ws.Range("R56").Value = 6736509
$ws.Range("S56").Value = 1
$ws.Range("T56").Value = "Dalarna"
$ws.Range("U56").Value = "Falun"
$ws.Range("V56").Value = "Dalarna"
$ws.Range("W56").Value = "Svärdsjö"
$ws.Range("Y56").Value = "2026-01-31"
$ws.Range("Z56").Value = ""
$ws.Range("AA56").Value = "2026-01-31"
$ws.Range("AB56").Value = ""
$ws.Range("AC56").Value = ""
$ws.Range("AD56").Value = $false
$ws.Range("AE56").Value = $false
$ws.Range("AF56").Value = ""
$ws.Range("AG56").Value = $false
$ws.Range("AT56").Value = ""
$ws.Range("AW56").Value = "Erik Danielsson"
$ws.Range("AX56").Value = "Erik Danielsson"
$ws.Range("AY56").Value = ""

# Row 57 <- content from original row 56
$ws.Range("A57").Value = 130979100
$ws.Range("B57").Value = 79244
$ws.Range("D57").Value = "NT"
$ws.Range("E57").Value = 6425
$ws.Range("F57").Value = "Garnlav"
$ws.Range("G57").Value = "Alectoria sarmentosa"
$ws.Range("H57").Value = "(Ach.) Ach."
$ws.Range("I57").Value = ""
$ws.Range("J57").Value = ""
$ws.Range("K57").Value = ""
$ws.Range("L57").Value = ""
$ws.Range("M57").Value = ""
$ws.Range("N57").Value = ""
$ws.Range("P57").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q57").Value = 571473
$ws.Range("R57").Value = 6736490
$ws.Range("S57").Value = 1
$ws.Range("T57").Value = "Dalarna"
$ws.Range("U57").Value = "Falun"
$ws.Range("V57").Value = "Dalarna"
$ws.Range("W57").Value = "Svärdsjö"
$ws.Range("Y57").Value = "2026-01-31"
$ws.Range("Z57").Value = ""
$ws.Range("AA57").Value = "2026-01-31"
$ws.Range("AB57").Value = ""
$ws.Range("AC57").Value = ""
$ws.Range("AD57").Value = $false
$ws.Range("AE57").Value = $false
$ws.Range("AF57").Value = ""
$ws.Range("AG57").Value = $false
$ws.Range("AT57").Value = ""
$ws.Range("AW57").Value = "Erik Danielsson"
$ws.Range("AX57").Value = "Erik Danielsson"
$ws.Range("AY57").Value = ""

# Row 58 <- content from original row 57
$ws.Range("A58").Value = 130979104
$ws.Range("B58").Value = 79244
$ws.Range("D58").Value = "NT"
$ws.Range("E58").Value = 6425
$ws.Range("F58").Value = "Garnlav"
$ws.Range("G58").Value = "Alectoria sarmentosa"
$ws.Range("H58").Value = "(Ach.) Ach."
$ws.Range("I58").Value = ""
$ws.Range("J58").Value = ""
$ws.Range("K58").Value = ""
$ws.Range("L58").Value = ""
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = ""
$ws.Range("P58").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q58").Value = 571129
$ws.Range("R58").Value = 6736573
$ws.Range("S58").Value = 1
$ws.Range("T58").Value = "Dalarna"
$ws.Range("U58").Value = "Falun"
$ws.Range("V58").Value = "Dalarna"
$ws.Range("W58").Value = "Svärdsjö"
$ws.Range("Y58").Value = "2026-01-31"
$ws.Range("Z58").Value = ""
$ws.Range("AA58").Value = "2026-01-31"
$ws.Range("AB58").Value = ""
$ws.Range("AC58").Value = ""
$ws.Range("AD58").Value = $false
$ws.Range("AE58").Value = $false
$ws.Range("AF58").Value = ""
$ws.Range("AG58").Value = $false
$ws.Range("AT58").Value = ""
$ws.Range("AW58").Value = "Erik Danielsson"
$ws.Range("AX58").Value = "Erik Danielsson"
$ws.Range("AY58").Value = ""

# Row 59 <- content from original row 58
$ws.Range("A59").Value = 130979094
$ws.Range("B59").Value = 79244
$ws.Range("D59").Value = "NT"
$ws.Range("E59").Value = 6425
$ws.Range("F59").Value = "Garnlav"
$ws.Range("G59").Value = "Alectoria sarmentosa"
$ws.Range("H59").Value = "(Ach.) Ach."
$ws.Range("I59").Value = ""
$ws.Range("J59").Value = ""
$ws.Range("K59").Value = ""
$ws.Range("L59").Value = ""
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = ""
$ws.Range("P59").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q59").Value = 571278
$ws.Range("R59").Value = 6736783
$ws.Range("S59").Value = 1
$ws.Range("T59").Value = "Dalarna"
$ws.Range("U59").Value = "Falun"
$ws.Range("V59").Value = "Dalarna"
$ws.Range("W59").Value = "Svärdsjö"
$ws.Range("Y59").Value = "2026-01-31"
$ws.Range("Z59").Value = ""
$ws.Range("AA59").Value = "2026-01-31"
$ws.Range("AB59").Value = ""
$ws.Range("AC59").Value = ""
$ws.Range("AD59").Value = $false
$ws.Range("AE59").Value = $false
$ws.Range("AF59").Value = ""
$ws.Range("AG59").Value = $false
$ws.Range("AT59").Value = ""
$ws.Range("AW59").Value = "Erik Danielsson"
$ws.Range("AX59").Value = "Erik Danielsson"
$ws.Range("AY59").Value = ""

